$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update F column "want to go" counts
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 577
$ws1.Range("F4").Value = 550
$ws1.Range("F5").Value = 43
$ws1.Range("F7").Value = 54
$ws1.Range("F9").Value = 52
$ws1.Range("F11").Value = 4692
$ws1.Range("F12").Value = 4484
$ws1.Range("F13").Value = 16
$ws1.Range("F15").Value = 25

# Sheet "全部类型" (sheet4): update F column "want to go" counts
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 577
$ws4.Range("F4").Value = 550
$ws4.Range("F5").Value = 43
$ws4.Range("F7").Value = 54
$ws4.Range("F9").Value = 52
$ws4.Range("F11").Value = 4692
$ws4.Range("F12").Value = 4484
$ws4.Range("F13").Value = 16
$ws4.Range("F15").Value = 25

$wb.Save()
